# Auto-generated edit script: applies market-data value corrections
# to the Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 949.2286
$ws.Cells.Item(15, 9).Value = 949.2286
$ws.Cells.Item(15, 11).Value = 2847.6858
$ws.Cells.Item(15, 13).Value = -2678.6858
$ws.Cells.Item(17, 8).Value = 1714.4231
$ws.Cells.Item(17, 10).Value = 1714.4231
$ws.Cells.Item(17, 12).Value = 5143.2693
$ws.Cells.Item(17, 14).Value = -5479.2693
$ws.Cells.Item(112, 8).Value = 2021.8889
$ws.Cells.Item(112, 10).Value = 2242.4285
$ws.Cells.Item(112, 12).Value = 6727.2855
$ws.Cells.Item(112, 14).Value = -8943.2855
$ws.Cells.Item(125, 8).Value = 1808
$ws.Cells.Item(125, 10).Value = 989.75
$ws.Cells.Item(125, 12).Value = 8907.75
$ws.Cells.Item(125, 14).Value = -13827.75
$ws.Cells.Item(132, 8).Value = 223101.28
$ws.Cells.Item(132, 9).Value = 909.675
$ws.Cells.Item(132, 11).Value = 2729.025
$ws.Cells.Item(132, 13).Value = -199.0249999999996
$ws.Cells.Item(138, 8).Value = 2614.2327
$ws.Cells.Item(138, 9).Value = 2531.4443
$ws.Cells.Item(138, 10).Value = 2673.84
$ws.Cells.Item(138, 11).Value = 7594.3329
$ws.Cells.Item(138, 12).Value = 8021.52
$ws.Cells.Item(138, 13).Value = -2454.3329
$ws.Cells.Item(138, 14).Value = -18301.52

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2178.4
$ws.Cells.Item(32, 9).Value = 1172.1285
$ws.Cells.Item(32, 11).Value = 1172.1285
$ws.Cells.Item(32, 13).Value = -885.1285
$ws.Cells.Item(45, 8).Value = 1809.7142
$ws.Cells.Item(45, 9).Value = 1869.75
$ws.Cells.Item(45, 11).Value = 1869.75
$ws.Cells.Item(45, 13).Value = -1492.75
$ws.Cells.Item(61, 8).Value = 3033.2693
$ws.Cells.Item(61, 9).Value = 2930.72
$ws.Cells.Item(61, 10).Value = 5597
$ws.Cells.Item(61, 11).Value = 2930.72
$ws.Cells.Item(61, 12).Value = 5597
$ws.Cells.Item(61, 13).Value = -2718.72
$ws.Cells.Item(61, 14).Value = -6021
$ws.Cells.Item(122, 8).Value = 3092.4546
$ws.Cells.Item(122, 9).Value = 3101.75
$ws.Cells.Item(122, 11).Value = 9305.25
$ws.Cells.Item(122, 13).Value = -6855.25
$ws.Cells.Item(136, 8).Value = 3033.2693
$ws.Cells.Item(136, 9).Value = 2930.72
$ws.Cells.Item(136, 11).Value = 8792.16
$ws.Cells.Item(136, 13).Value = -6242.16
$ws.Cells.Item(136, 14).Value = -21891

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 4478.353
$ws.Cells.Item(20, 9).Value = 2858.7856
$ws.Cells.Item(20, 11).Value = 2858.7856
$ws.Cells.Item(20, 13).Value = -2611.7856
$ws.Cells.Item(64, 8).Value = 1432
$ws.Cells.Item(64, 10).Value = 900
$ws.Cells.Item(64, 12).Value = 900
$ws.Cells.Item(64, 14).Value = -1350
$ws.Cells.Item(67, 8).Value = 1432
$ws.Cells.Item(67, 10).Value = 900
$ws.Cells.Item(67, 12).Value = 900
$ws.Cells.Item(67, 14).Value = -2460
$ws.Cells.Item(80, 8).Value = 2157.5454
$ws.Cells.Item(80, 9).Value = 510.25
$ws.Cells.Item(80, 11).Value = 510.25
$ws.Cells.Item(80, 13).Value = 487.75
$ws.Cells.Item(81, 8).Value = 44894.5
$ws.Cells.Item(81, 10).Value = 44894.5
$ws.Cells.Item(81, 12).Value = 44894.5
$ws.Cells.Item(81, 14).Value = -47016.5
$ws.Cells.Item(83, 8).Value = 2157.5454
$ws.Cells.Item(83, 9).Value = 510.25
$ws.Cells.Item(83, 11).Value = 2551.25
$ws.Cells.Item(83, 13).Value = 2440.75
$ws.Cells.Item(84, 8).Value = 44894.5
$ws.Cells.Item(84, 10).Value = 44894.5
$ws.Cells.Item(84, 12).Value = 134683.5
$ws.Cells.Item(84, 14).Value = -145291.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 77057.46000000001
$ws.Cells.Item(7, 9).Value = 111216.22
$ws.Cells.Item(7, 11).Value = 111216.22
$ws.Cells.Item(7, 13).Value = -111103.22
$ws.Cells.Item(22, 8).Value = 599.6429000000001
$ws.Cells.Item(22, 10).Value = 1600
$ws.Cells.Item(22, 12).Value = 1600
$ws.Cells.Item(22, 14).Value = -2300
$ws.Cells.Item(31, 8).Value = 1702
$ws.Cells.Item(31, 9).Value = 1702
$ws.Cells.Item(31, 11).Value = 1702
$ws.Cells.Item(31, 13).Value = -1407
$ws.Cells.Item(34, 8).Value = 1702
$ws.Cells.Item(34, 9).Value = 1702
$ws.Cells.Item(34, 11).Value = 1702
$ws.Cells.Item(34, 13).Value = -1500
$ws.Cells.Item(68, 8).Value = 63670.57
$ws.Cells.Item(68, 10).Value = 63670.57
$ws.Cells.Item(68, 12).Value = 63670.57
$ws.Cells.Item(68, 14).Value = -65168.57
$ws.Cells.Item(71, 8).Value = 63670.57
$ws.Cells.Item(71, 10).Value = 63670.57
$ws.Cells.Item(71, 12).Value = 191011.71
$ws.Cells.Item(71, 14).Value = -198499.71
$ws.Cells.Item(105, 8).Value = 1954.7142
$ws.Cells.Item(105, 9).Value = 1819.2222
$ws.Cells.Item(105, 11).Value = 1819.2222
$ws.Cells.Item(105, 13).Value = -72.22219999999993
$ws.Cells.Item(134, 8).Value = 1703.425
$ws.Cells.Item(134, 9).Value = 1763
$ws.Cells.Item(134, 10).Value = 1365.8334
$ws.Cells.Item(134, 11).Value = 5289
$ws.Cells.Item(134, 12).Value = 4097.5002
$ws.Cells.Item(134, 13).Value = -2754
$ws.Cells.Item(134, 14).Value = -9167.5002

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(55, 8).Value = 2333.3333
$ws.Cells.Item(55, 10).Value = 2600
$ws.Cells.Item(55, 12).Value = 7800
$ws.Cells.Item(55, 14).Value = -8154
$ws.Cells.Item(92, 8).Value = 271.6
$ws.Cells.Item(92, 10).Value = 284.33334
$ws.Cells.Item(92, 12).Value = 853.0000200000001
$ws.Cells.Item(92, 14).Value = -3349.00002
$ws.Cells.Item(122, 8).Value = 958.2
$ws.Cells.Item(122, 9).Value = 673.25
$ws.Cells.Item(122, 10).Value = 1148.1666
$ws.Cells.Item(122, 11).Value = 6059.25
$ws.Cells.Item(122, 12).Value = 10333.4994
$ws.Cells.Item(122, 13).Value = -3609.25
$ws.Cells.Item(122, 14).Value = -15233.4994

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 12347.826
$ws.Cells.Item(80, 9).Value = 5758.375
$ws.Cells.Item(80, 11).Value = 5758.375
$ws.Cells.Item(80, 13).Value = -4760.375
$ws.Cells.Item(83, 8).Value = 12347.826
$ws.Cells.Item(83, 9).Value = 5758.375
$ws.Cells.Item(83, 11).Value = 28791.875
$ws.Cells.Item(83, 13).Value = -23799.875
$ws.Cells.Item(102, 8).Value = 3141.5
$ws.Cells.Item(102, 9).Value = 2078.6
$ws.Cells.Item(102, 11).Value = 2078.6
$ws.Cells.Item(102, 13).Value = -456.5999999999999
$ws.Cells.Item(132, 8).Value = 3025.5
$ws.Cells.Item(132, 9).Value = 1878
$ws.Cells.Item(132, 10).Value = 4173
$ws.Cells.Item(132, 11).Value = 5634
$ws.Cells.Item(132, 12).Value = 12519
$ws.Cells.Item(132, 13).Value = -3104
$ws.Cells.Item(132, 14).Value = -17579

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(21, 8).Value = 15004.667
$ws.Cells.Item(21, 9).Value = 0
$ws.Cells.Item(21, 10).Value = 15004.667
$ws.Cells.Item(21, 11).Value = 0
$ws.Cells.Item(21, 12).Value = 15004.667
$ws.Cells.Item(21, 13).ClearContents()
$ws.Cells.Item(21, 14).Value = -15352.667
$ws.Cells.Item(46, 8).Value = 2515.6924
$ws.Cells.Item(46, 9).Value = 1639.8
$ws.Cells.Item(46, 11).Value = 1639.8
$ws.Cells.Item(46, 13).Value = -1451.8
$ws.Cells.Item(122, 8).Value = 6268.4
$ws.Cells.Item(122, 9).Value = 6268.4
$ws.Cells.Item(122, 11).Value = 18805.2
$ws.Cells.Item(122, 13).Value = -16355.2
$ws.Cells.Item(132, 8).Value = 1996.8636
$ws.Cells.Item(132, 9).Value = 1856.55
$ws.Cells.Item(132, 11).Value = 5569.65
$ws.Cells.Item(132, 13).Value = -3039.65
$ws.Cells.Item(136, 8).Value = 2162.5356
$ws.Cells.Item(136, 9).Value = 2241.8262
$ws.Cells.Item(136, 10).Value = 1797.8
$ws.Cells.Item(136, 11).Value = 6725.4786
$ws.Cells.Item(136, 12).Value = 5393.4
$ws.Cells.Item(136, 13).Value = -4175.4786
$ws.Cells.Item(136, 14).Value = -10493.4

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(6, 8).Value = 5
$ws.Cells.Item(6, 9).Value = 5
$ws.Cells.Item(6, 11).Value = 5
$ws.Cells.Item(6, 13).Value = 110
$ws.Cells.Item(11, 8).Value = 3003
$ws.Cells.Item(11, 9).Value = 3003
$ws.Cells.Item(11, 10).Value = 0
$ws.Cells.Item(11, 11).Value = 3003
$ws.Cells.Item(11, 12).Value = 0
$ws.Cells.Item(11, 13).Value = -2861
$ws.Cells.Item(11, 14).ClearContents()
$ws.Cells.Item(23, 8).Value = 8419.799999999999
$ws.Cells.Item(23, 9).Value = 2033.3334
$ws.Cells.Item(23, 10).Value = 17999.5
$ws.Cells.Item(23, 11).Value = 2033.3334
$ws.Cells.Item(23, 12).Value = 17999.5
$ws.Cells.Item(23, 13).Value = -1804.3334
$ws.Cells.Item(23, 14).Value = -18457.5
$ws.Cells.Item(58, 8).Value = 40000
$ws.Cells.Item(58, 9).Value = 40000
$ws.Cells.Item(58, 11).Value = 40000
$ws.Cells.Item(58, 13).Value = -39692
$ws.Cells.Item(122, 8).Value = 2554.875
$ws.Cells.Item(122, 9).Value = 2097.4814
$ws.Cells.Item(122, 10).Value = 3504.8462
$ws.Cells.Item(122, 11).Value = 6292.4442
$ws.Cells.Item(122, 12).Value = 10514.5386
$ws.Cells.Item(122, 13).Value = -3842.4442
$ws.Cells.Item(122, 14).Value = -15414.5386
$ws.Cells.Item(132, 8).Value = 1364
$ws.Cells.Item(132, 9).Value = 1364
$ws.Cells.Item(132, 11).Value = 4092
$ws.Cells.Item(132, 13).Value = -1562
